$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: MuSCs -> ECs (sending cluster), target cluster -> ECs, TPM values updated ---
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Adm2"
$ws.Range("C2").Value = "Ramp2"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.072366
$ws.Range("H2").Value = 0.217098
$ws.Range("I2").Value = 1
$ws.Range("J2").Value = 1
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 82.374849
$ws.Range("N2").Value = 247.124547
$ws.Range("O2").Value = 0.6248537741508955
$ws.Range("P2").Value = 0.6248537741508954
$ws.Range("Q2").Value = 5.961138322734
$ws.Range("R2").Value = 53.65024490460601
$ws.Range("S2").Value = 0.6248537741508955
$ws.Range("T2").Value = 0.6248537741508954

# --- Row 3: MuSCs -> ECs (sending cluster), target cluster -> FAPs, TPM values updated ---
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Adm2"
$ws.Range("C3").Value = "Ramp2"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.072366
$ws.Range("H3").Value = 0.217098
$ws.Range("I3").Value = 1
$ws.Range("J3").Value = 1
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 32.15646633333333
$ws.Range("N3").Value = 96.46939900000001
$ws.Range("O3").Value = 0.2439226243891451
$ws.Range("P3").Value = 0.2439226243891451
$ws.Range("Q3").Value = 2.327034842678
$ws.Range("R3").Value = 20.943313584102
$ws.Range("S3").Value = 0.2439226243891451
$ws.Range("T3").Value = 0.2439226243891451

# --- Row 4: MuSCs -> ECs (sending cluster), target cluster -> MuSCs, TPM values updated ---
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Adm2"
$ws.Range("C4").Value = "Ramp2"
$ws.Range("D4").Value = "MuSCs"
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.072366
$ws.Range("H4").Value = 0.217098
$ws.Range("I4").Value = 1
$ws.Range("J4").Value = 1
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 17.22395566666667
$ws.Range("N4").Value = 51.671867
$ws.Range("O4").Value = 0.130652181275918
$ws.Range("P4").Value = 0.1306521812759179
$ws.Range("Q4").Value = 1.246428775774
$ws.Range("R4").Value = 11.217858981966
$ws.Range("S4").Value = 0.130652181275918
$ws.Range("T4").Value = 0.1306521812759179

# --- Row 5 (new): ECs -> Adm2/Ramp2 -> Resolving-Mac ---
$ws.Range("A5").Value = "ECs"
$ws.Range("B5").Value = "Adm2"
$ws.Range("C5").Value = "Ramp2"
$ws.Range("D5").Value = "Resolving-Mac"
$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 0.3333333333333333
$ws.Range("G5").Value = 0.072366
$ws.Range("H5").Value = 0.217098
$ws.Range("I5").Value = 1
$ws.Range("J5").Value = 1
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 0.07533066666666667
$ws.Range("N5").Value = 0.225992
$ws.Range("O5").Value = 0.0005714201840414873
$ws.Range("P5").Value = 0.0005714201840414872
$ws.Range("Q5").Value = 0.005451379024
$ws.Range("R5").Value = 0.049062411216
$ws.Range("S5").Value = 0.0005714201840414873
$ws.Range("T5").Value = 0.0005714201840414872
